$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2899.5652
$ws.Range("I137").Value = 2311.182
$ws.Range("J137").Value = 3438.9167
$ws.Range("K137").Value = 6933.545999999999
$ws.Range("L137").Value = 10316.7501
$ws.Range("M137").Value = -4383.545999999999
$ws.Range("N137").Value = -15416.7501

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7886.7144
$ws.Range("I61").Value = 7880.2
$ws.Range("J61").Value = 7890.3335
$ws.Range("K61").Value = 7880.2
$ws.Range("L61").Value = 7890.3335
$ws.Range("M61").Value = -7668.2
$ws.Range("N61").Value = -8314.333500000001
$ws.Range("H64").Value = 33888.11
$ws.Range("J64").Value = 33888.11
$ws.Range("L64").Value = 33888.11
$ws.Range("N64").Value = -34384.11
$ws.Range("H67").Value = 33888.11
$ws.Range("J67").Value = 33888.11
$ws.Range("L67").Value = 33888.11
$ws.Range("N67").Value = -35604.11
$ws.Range("H74").Value = 1445.9474
$ws.Range("I74").Value = 1335.875
$ws.Range("J74").Value = 2033
$ws.Range("K74").Value = 1335.875
$ws.Range("L74").Value = 2033
$ws.Range("M74").Value = -461.875
$ws.Range("N74").Value = -3781
$ws.Range("H77").Value = 1445.9474
$ws.Range("I77").Value = 1335.875
$ws.Range("J77").Value = 2033
$ws.Range("K77").Value = 6679.375
$ws.Range("L77").Value = 10165
$ws.Range("M77").Value = -2311.375
$ws.Range("N77").Value = -18901
$ws.Range("H97").Value = 1145.1666
$ws.Range("I97").Value = 1145.1666
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1145.1666
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -649.1666
$ws.Range("H122").Value = 5672.1763
$ws.Range("I122").Value = 4643.1
$ws.Range("J122").Value = 7142.2856
$ws.Range("K122").Value = 13929.3
$ws.Range("L122").Value = 21426.8568
$ws.Range("M122").Value = -11479.3
$ws.Range("N122").Value = -26326.8568
$ws.Range("H136").Value = 7886.7144
$ws.Range("I136").Value = 7880.2
$ws.Range("J136").Value = 7890.3335
$ws.Range("K136").Value = 23640.6
$ws.Range("L136").Value = 23671.0005
$ws.Range("M136").Value = -21090.6
$ws.Range("N136").Value = -28771.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6474.5625
$ws.Range("I86").Value = 4715.5
$ws.Range("J86").Value = 11751.75
$ws.Range("K86").Value = 4715.5
$ws.Range("L86").Value = 11751.75
$ws.Range("M86").Value = -3592.5
$ws.Range("N86").Value = -13997.75
$ws.Range("H89").Value = 6474.5625
$ws.Range("I89").Value = 4715.5
$ws.Range("J89").Value = 11751.75
$ws.Range("K89").Value = 23577.5
$ws.Range("L89").Value = 58758.75
$ws.Range("M89").Value = -17961.5
$ws.Range("N89").Value = -69990.75
$ws.Range("H134").Value = 2908.861
$ws.Range("I134").Value = 1570.5217
$ws.Range("J134").Value = 5276.6924
$ws.Range("K134").Value = 4711.5651
$ws.Range("L134").Value = 15830.0772
$ws.Range("M134").Value = -2176.5651
$ws.Range("N134").Value = -20900.0772

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3437.6072
$ws.Range("I31").Value = 2433.8823
$ws.Range("J31").Value = 4988.8184
$ws.Range("K31").Value = 2433.8823
$ws.Range("L31").Value = 4988.8184
$ws.Range("M31").Value = -2138.8823
$ws.Range("N31").Value = -5578.8184
$ws.Range("H34").Value = 3437.6072
$ws.Range("I34").Value = 2433.8823
$ws.Range("J34").Value = 4988.8184
$ws.Range("K34").Value = 2433.8823
$ws.Range("L34").Value = 4988.8184
$ws.Range("M34").Value = -2231.8823
$ws.Range("N34").Value = -5392.8184
$ws.Range("H99").Value = 16035698
$ws.Range("I99").Value = 3057782
$ws.Range("J99").Value = 33339586
$ws.Range("K99").Value = 3057782
$ws.Range("L99").Value = 33339586
$ws.Range("M99").Value = -3056284
$ws.Range("N99").Value = -33342582
$ws.Range("H122").Value = 6555556.5
$ws.Range("I122").Value = 44219630
$ws.Range("K122").Value = 132658890
$ws.Range("M122").Value = -132656440
$ws.Range("H126").Value = 16035698
$ws.Range("I126").Value = 3057782
$ws.Range("J126").Value = 33339586
$ws.Range("K126").Value = 9173346
$ws.Range("L126").Value = 100018758
$ws.Range("M126").Value = -9170876
$ws.Range("N126").Value = -100023698
$ws.Range("H132").Value = 3980.8
$ws.Range("I132").Value = 3847.9546
$ws.Range("K132").Value = 11543.8638
$ws.Range("M132").Value = -9013.863799999999
$ws.Range("H134").Value = 2265.457
$ws.Range("I134").Value = 1344.4828
$ws.Range("J134").Value = 6716.8335
$ws.Range("K134").Value = 4033.4484
$ws.Range("L134").Value = 20150.5005
$ws.Range("M134").Value = -1498.4484
$ws.Range("N134").Value = -25220.5005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3007644.5
$ws.Range("J7").Value = 3007644.5
$ws.Range("L7").Value = 3007644.5
$ws.Range("N7").Value = -3007868.5
$ws.Range("H8").Value = 3007644.5
$ws.Range("J8").Value = 3007644.5
$ws.Range("L8").Value = 3007644.5
$ws.Range("N8").Value = -3007922.5
$ws.Range("H122").Value = 5643.5454
$ws.Range("I122").Value = 1680
$ws.Range("J122").Value = 7908.4287
$ws.Range("K122").Value = 5040
$ws.Range("L122").Value = 23725.2861
$ws.Range("M122").Value = -2590
$ws.Range("N122").Value = -28625.2861
$ws.Range("H126").Value = 6381.6665
$ws.Range("I126").Value = 5597.75
$ws.Range("J126").Value = 6773.625
$ws.Range("K126").Value = 16793.25
$ws.Range("L126").Value = 20320.875
$ws.Range("M126").Value = -14323.25
$ws.Range("N126").Value = -25260.875
$ws.Range("H132").Value = 4096.648
$ws.Range("I132").Value = 2999.4707
$ws.Range("J132").Value = 5961.85
$ws.Range("K132").Value = 8998.4121
$ws.Range("L132").Value = 17885.55
$ws.Range("M132").Value = -6468.4121
$ws.Range("N132").Value = -22945.55

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 999999
$ws.Range("J3").Value = 999999
$ws.Range("L3").Value = 999999
$ws.Range("N3").Value = -1000223
$ws.Range("H5").Value = 22600
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H15").Value = 999999
$ws.Range("J15").Value = 999999
$ws.Range("L15").Value = 999999
$ws.Range("N15").Value = -1000339

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 18400
$ws.Range("I21").Value = 14000
$ws.Range("J21").Value = 22800
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 22800
$ws.Range("M21").Value = -13765
$ws.Range("N21").Value = -23270
$ws.Range("H35").Value = 18400
$ws.Range("I35").Value = 14000
$ws.Range("J35").Value = 22800
$ws.Range("K35").Value = 14000
$ws.Range("L35").Value = 22800
$ws.Range("M35").Value = -13710
$ws.Range("N35").Value = -23380
$ws.Range("H122").Value = 3278
$ws.Range("I122").Value = 2650.074
$ws.Range("J122").Value = 5700
$ws.Range("K122").Value = 7950.222
$ws.Range("L122").Value = 17100
$ws.Range("M122").Value = -5500.222
$ws.Range("N122").Value = -22000
$ws.Range("H126").Value = 1918.75
$ws.Range("I126").Value = 1780.5
$ws.Range("J126").Value = 2195.25
$ws.Range("K126").Value = 5341.5
$ws.Range("L126").Value = 6585.75
$ws.Range("M126").Value = -2871.5
$ws.Range("N126").Value = -11525.75
$ws.Range("H132").Value = 1929.7391
$ws.Range("I132").Value = 968.7941
$ws.Range("J132").Value = 4652.4165
$ws.Range("K132").Value = 2906.3823
$ws.Range("L132").Value = 13957.2495
$ws.Range("M132").Value = -376.3822999999998
$ws.Range("N132").Value = -19017.2495
$ws.Range("H136").Value = 3113.9019
$ws.Range("I136").Value = 2438.4285
$ws.Range("J136").Value = 6266.1113
$ws.Range("K136").Value = 7315.2855
$ws.Range("L136").Value = 18798.3339
$ws.Range("M136").Value = -4765.2855
$ws.Range("N136").Value = -23898.3339

Write-Output "Applied all cell updates."